$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Adding new releases (rows) and updating analysis metrics for existing releases.
# Final data block spans A2:M10 (dimension grows from A1:M8 to A1:M10).

# Row 2
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = 0.4
$ws.Cells.Item(2, 3).Value = 1
$ws.Cells.Item(2, 4).Value = 1
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = "Archives"
$ws.Cells.Item(2, 7).Value = "03-10-2021-15-47"
$ws.Cells.Item(2, 8).Value = 200
$ws.Cells.Item(2, 9).Value = 0.792
$ws.Cells.Item(2, 10).Value = 1
$ws.Cells.Item(2, 11).Value = 0.396
$ws.Cells.Item(2, 12).Value = 0.5
$ws.Cells.Item(2, 13).Value = 0.896

# Row 3
$ws.Cells.Item(3, 1).Value = 0.4285714285714285
$ws.Cells.Item(3, 2).Value = 0.5714285714285714
$ws.Cells.Item(3, 3).Value = 1
$ws.Cells.Item(3, 4).Value = 1
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = "Archives"
$ws.Cells.Item(3, 7).Value = "03-11-2021-02-10"
$ws.Cells.Item(3, 8).Value = 428
$ws.Cells.Item(3, 9).Value = 0.6599999999999999
$ws.Cells.Item(3, 10).Value = 1
$ws.Cells.Item(3, 11).Value = 0.33
$ws.Cells.Item(3, 12).Value = 0.5
$ws.Cells.Item(3, 13).Value = 0.83

# Row 4
$ws.Cells.Item(4, 1).Value = 1
$ws.Cells.Item(4, 2).Value = 0.4
$ws.Cells.Item(4, 3).Value = 1
$ws.Cells.Item(4, 4).Value = 1
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = "Archives"
$ws.Cells.Item(4, 7).Value = "13-09-2021-15-00"
$ws.Cells.Item(4, 8).Value = 109
$ws.Cells.Item(4, 9).Value = 0.792
$ws.Cells.Item(4, 10).Value = 1
$ws.Cells.Item(4, 11).Value = 0.396
$ws.Cells.Item(4, 12).Value = 0.5
$ws.Cells.Item(4, 13).Value = 0.896

# Row 5
$ws.Cells.Item(5, 1).Value = 1
$ws.Cells.Item(5, 2).Value = 0.4
$ws.Cells.Item(5, 3).Value = 1
$ws.Cells.Item(5, 4).Value = 1
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = "Profile"
$ws.Cells.Item(5, 7).Value = "03-10-2021-15-48"
$ws.Cells.Item(5, 8).Value = 112
$ws.Cells.Item(5, 9).Value = 0.792
$ws.Cells.Item(5, 10).Value = 1
$ws.Cells.Item(5, 11).Value = 0.396
$ws.Cells.Item(5, 12).Value = 0.5
$ws.Cells.Item(5, 13).Value = 0.896

# Row 6
$ws.Cells.Item(6, 1).Value = 1
$ws.Cells.Item(6, 2).Value = 0.4
$ws.Cells.Item(6, 3).Value = 1
$ws.Cells.Item(6, 4).Value = 1
$ws.Cells.Item(6, 5).Value = 1
$ws.Cells.Item(6, 6).Value = "Profile"
$ws.Cells.Item(6, 7).Value = "13-09-2021-14-00"
$ws.Cells.Item(6, 8).Value = 112
$ws.Cells.Item(6, 9).Value = 0.792
$ws.Cells.Item(6, 10).Value = 1
$ws.Cells.Item(6, 11).Value = 0.396
$ws.Cells.Item(6, 12).Value = 0.5
$ws.Cells.Item(6, 13).Value = 0.896

# Row 7
$ws.Cells.Item(7, 1).Value = 0.9090909090909091
$ws.Cells.Item(7, 2).Value = 0.04545454545454546
$ws.Cells.Item(7, 3).Value = 1
$ws.Cells.Item(7, 4).Value = 1
$ws.Cells.Item(7, 5).Value = 1
$ws.Cells.Item(7, 6).Value = "Frontend"
$ws.Cells.Item(7, 7).Value = "03-10-2021-15-49"
$ws.Cells.Item(7, 8).Value = 1401
$ws.Cells.Item(7, 9).Value = 0.645
$ws.Cells.Item(7, 10).Value = 1
$ws.Cells.Item(7, 11).Value = 0.3225
$ws.Cells.Item(7, 12).Value = 0.5
$ws.Cells.Item(7, 13).Value = 0.8225

# Row 8
$ws.Cells.Item(8, 1).Value = 0.9591836734693877
$ws.Cells.Item(8, 2).Value = 0.02040816326530612
$ws.Cells.Item(8, 3).Value = 0.9183673469387755
$ws.Cells.Item(8, 4).Value = 1
$ws.Cells.Item(8, 5).Value = 0.9387755102040817
$ws.Cells.Item(8, 6).Value = "Frontend"
$ws.Cells.Item(8, 7).Value = "03-11-2021-02-08"
$ws.Cells.Item(8, 8).Value = 4888
$ws.Cells.Item(8, 9).Value = 0.626326530612245
$ws.Cells.Item(8, 10).Value = 0.9571428571428571
$ws.Cells.Item(8, 11).Value = 0.3131632653061225
$ws.Cells.Item(8, 12).Value = 0.4785714285714285
$ws.Cells.Item(8, 13).Value = 0.791734693877551

# Row 9
$ws.Cells.Item(9, 1).Value = 0.5
$ws.Cells.Item(9, 2).Value = 0.5
$ws.Cells.Item(9, 3).Value = 1
$ws.Cells.Item(9, 4).Value = 1
$ws.Cells.Item(9, 5).Value = 1
$ws.Cells.Item(9, 6).Value = "Frontend"
$ws.Cells.Item(9, 7).Value = "13-09-2021-20-00"
$ws.Cells.Item(9, 8).Value = 35
$ws.Cells.Item(9, 9).Value = 0.66
$ws.Cells.Item(9, 10).Value = 1
$ws.Cells.Item(9, 11).Value = 0.33
$ws.Cells.Item(9, 12).Value = 0.5
$ws.Cells.Item(9, 13).Value = 0.8300000000000001

# Row 10
$ws.Cells.Item(10, 1).Value = 0.9591836734693877
$ws.Cells.Item(10, 2).Value = 0.02040816326530612
$ws.Cells.Item(10, 3).Value = 0.9183673469387755
$ws.Cells.Item(10, 4).Value = 1
$ws.Cells.Item(10, 5).Value = 0.9591836734693877
$ws.Cells.Item(10, 6).Value = "Frontend"
$ws.Cells.Item(10, 7).Value = "17-10-2021-15-30"
$ws.Cells.Item(10, 8).Value = 4858
$ws.Cells.Item(10, 9).Value = 0.626326530612245
$ws.Cells.Item(10, 10).Value = 0.9714285714285713
$ws.Cells.Item(10, 11).Value = 0.3131632653061225
$ws.Cells.Item(10, 12).Value = 0.4857142857142857
$ws.Cells.Item(10, 13).Value = 0.7988775510204081
